$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:D2").WrapText = $true
$ws.Range("E2:E18").WrapText = $true

$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "RESIDENT-REG"
$ws.Range("C2").Value = "Resident services virtual registration machine"
$ws.Range("D2").Value = "Non-existent machine created for Resident services"
$ws.Range("E2").Value = "'TRUE"

$ws.Rows.Item(2).RowHeight = 22.8

$ws.Columns.Item(1).ColumnWidth = 9.48
$ws.Columns.Item(2).ColumnWidth = 16.77
$ws.Columns.Item(3).ColumnWidth = 24.97
$ws.Columns.Item(4).ColumnWidth = 23.87
$ws.Columns.Item(5).ColumnWidth = 8.19

$ws.Range("E2").Select()
